# Auto-generated Excel COM-interop script applying data updates to Lamia_Profits sheets
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$ws_ALC.Range("H19").Value = 1283.3334
$ws_ALC.Range("J19").Value = 1375
$ws_ALC.Range("L19").Value = 1375
$ws_ALC.Range("N19").Value = -1725
$ws_ALC.Range("H64").Value = 6101
$ws_ALC.Range("J64").Value = 0
$ws_ALC.Range("L64").Value = 0
$ws_ALC.Range("N64").ClearContents()
$ws_ALC.Range("H67").Value = 6101
$ws_ALC.Range("J67").Value = 0
$ws_ALC.Range("L67").Value = 0
$ws_ALC.Range("N67").ClearContents()
$ws_ALC.Range("H74").Value = 8046.722
$ws_ALC.Range("I74").Value = 6094.364
$ws_ALC.Range("K74").Value = 6094.364
$ws_ALC.Range("M74").Value = -5158.364
$ws_ALC.Range("H77").Value = 8046.722
$ws_ALC.Range("I77").Value = 6094.364
$ws_ALC.Range("K77").Value = 30471.82
$ws_ALC.Range("M77").Value = -25791.82
$ws_ALC.Range("H118").Value = 904.2727
$ws_ALC.Range("I118").Value = 390.5
$ws_ALC.Range("J118").Value = 1520.8
$ws_ALC.Range("K118").Value = 1171.5
$ws_ALC.Range("L118").Value = 4562.4
$ws_ALC.Range("M118").Value = 485.5
$ws_ALC.Range("N118").Value = -7876.4
$ws_ALC.Range("H138").Value = 2733.9656
$ws_ALC.Range("J138").Value = 3499.6667
$ws_ALC.Range("L138").Value = 10499.0001
$ws_ALC.Range("N138").Value = -20779.0001
# --- ARM ---
$ws_ARM.Range("H2").Value = 11669.667
$ws_ARM.Range("I2").Value = 628.375
$ws_ARM.Range("K2").Value = 628.375
$ws_ARM.Range("M2").Value = -515.375
$ws_ARM.Range("H32").Value = 22729548
$ws_ARM.Range("I32").Value = 23811210
$ws_ARM.Range("K32").Value = 23811210
$ws_ARM.Range("M32").Value = -23810923
$ws_ARM.Range("H97").Value = 1316.1428
$ws_ARM.Range("I97").Value = 1538.6364
$ws_ARM.Range("K97").Value = 1538.6364
$ws_ARM.Range("M97").Value = -1042.6364
$ws_ARM.Range("H116").Value = 11669.667
$ws_ARM.Range("I116").Value = 628.375
$ws_ARM.Range("K116").Value = 628.375
$ws_ARM.Range("M116").Value = 1665.625
# --- BSM ---
$ws_BSM.Range("H3").Value = 11669.667
$ws_BSM.Range("I3").Value = 628.375
$ws_BSM.Range("K3").Value = 628.375
$ws_BSM.Range("M3").Value = -514.375
$ws_BSM.Range("H20").Value = 2793.4546
$ws_BSM.Range("I20").Value = 1948.2858
$ws_BSM.Range("J20").Value = 4272.5
$ws_BSM.Range("K20").Value = 1948.2858
$ws_BSM.Range("L20").Value = 4272.5
$ws_BSM.Range("M20").Value = -1701.2858
$ws_BSM.Range("N20").Value = -4766.5
$ws_BSM.Range("H86").Value = 2597.75
$ws_BSM.Range("I86").Value = 2492.95
$ws_BSM.Range("K86").Value = 2492.95
$ws_BSM.Range("M86").Value = -1369.95
$ws_BSM.Range("H89").Value = 2597.75
$ws_BSM.Range("I89").Value = 2492.95
$ws_BSM.Range("K89").Value = 12464.75
$ws_BSM.Range("M89").Value = -6848.75
$ws_BSM.Range("H134").Value = 2835.5386
$ws_BSM.Range("I134").Value = 1335.2174
$ws_BSM.Range("K134").Value = 4005.6522
$ws_BSM.Range("M134").Value = -1470.6522
# --- CRP ---
$ws_CRP.Range("H132").Value = 3580.8
$ws_CRP.Range("I132").Value = 1774
$ws_CRP.Range("J132").Value = 4785.3335
$ws_CRP.Range("K132").Value = 5322
$ws_CRP.Range("L132").Value = 14356.0005
$ws_CRP.Range("M132").Value = -2792
$ws_CRP.Range("N132").Value = -19416.0005
$ws_CRP.Range("H134").Value = 3531.3125
$ws_CRP.Range("I134").Value = 2007.25
$ws_CRP.Range("J134").Value = 8103.5
$ws_CRP.Range("K134").Value = 6021.75
$ws_CRP.Range("L134").Value = 24310.5
$ws_CRP.Range("M134").Value = -3486.75
$ws_CRP.Range("N134").Value = -29380.5
# --- CUL ---
$ws_CUL.Range("H11").Value = 250262.5
$ws_CUL.Range("I11").Value = 333633.34
$ws_CUL.Range("K11").Value = 1000900.02
$ws_CUL.Range("M11").Value = -1000760.02
$ws_CUL.Range("H68").Value = 15626405
$ws_CUL.Range("I68").Value = 25001398
$ws_CUL.Range("J68").Value = 1415
$ws_CUL.Range("K68").Value = 75004194
$ws_CUL.Range("L68").Value = 4245
$ws_CUL.Range("M68").Value = -75003383
$ws_CUL.Range("N68").Value = -5867
$ws_CUL.Range("H70").Value = 11338
$ws_CUL.Range("I70").Value = 9500
$ws_CUL.Range("K70").Value = 28500
$ws_CUL.Range("M70").Value = -28185
$ws_CUL.Range("H71").Value = 15626405
$ws_CUL.Range("I71").Value = 25001398
$ws_CUL.Range("J71").Value = 1415
$ws_CUL.Range("K71").Value = 225012582
$ws_CUL.Range("L71").Value = 12735
$ws_CUL.Range("M71").Value = -225008526
$ws_CUL.Range("N71").Value = -20847
$ws_CUL.Range("H73").Value = 11338
$ws_CUL.Range("I73").Value = 9500
$ws_CUL.Range("K73").Value = 28500
$ws_CUL.Range("M73").Value = -27408
$ws_CUL.Range("H75").Value = 83339380
$ws_CUL.Range("I75").Value = 250000660
$ws_CUL.Range("J75").Value = 8744.875
$ws_CUL.Range("K75").Value = 750001980
$ws_CUL.Range("L75").Value = 26234.625
$ws_CUL.Range("M75").Value = -750000982
$ws_CUL.Range("N75").Value = -28230.625
$ws_CUL.Range("H78").Value = 83339380
$ws_CUL.Range("I78").Value = 250000660
$ws_CUL.Range("J78").Value = 8744.875
$ws_CUL.Range("K78").Value = 2250005940
$ws_CUL.Range("L78").Value = 78703.875
$ws_CUL.Range("M78").Value = -2250000948
$ws_CUL.Range("N78").Value = -88687.875
$ws_CUL.Range("H122").Value = 6741502
$ws_CUL.Range("I122").Value = 4273989.5
$ws_CUL.Range("J122").Value = 7520716
$ws_CUL.Range("K122").Value = 38465905.5
$ws_CUL.Range("L122").Value = 67686444
$ws_CUL.Range("M122").Value = -38463455.5
$ws_CUL.Range("N122").Value = -67691344
$ws_CUL.Range("H132").Value = 4076.9092
$ws_CUL.Range("I132").Value = 2563
$ws_CUL.Range("J132").Value = 5590.8184
$ws_CUL.Range("K132").Value = 23067
$ws_CUL.Range("L132").Value = 50317.3656
$ws_CUL.Range("M132").Value = -20537
$ws_CUL.Range("N132").Value = -55377.3656
# --- GSM ---
$ws_GSM.Range("H21").Value = 26926.25
$ws_GSM.Range("I21").Value = 25002
$ws_GSM.Range("K21").Value = 25002
$ws_GSM.Range("M21").Value = -24829
$ws_GSM.Range("H30").Value = 26926.25
$ws_GSM.Range("I30").Value = 25002
$ws_GSM.Range("K30").Value = 25002
$ws_GSM.Range("M30").Value = -24897
$ws_GSM.Range("H70").Value = 10290.385
$ws_GSM.Range("I70").Value = 7032.875
$ws_GSM.Range("J70").Value = 15502.4
$ws_GSM.Range("K70").Value = 7032.875
$ws_GSM.Range("L70").Value = 15502.4
$ws_GSM.Range("M70").Value = -6762.875
$ws_GSM.Range("N70").Value = -16042.4
$ws_GSM.Range("H73").Value = 10290.385
$ws_GSM.Range("I73").Value = 7032.875
$ws_GSM.Range("J73").Value = 15502.4
$ws_GSM.Range("K73").Value = 7032.875
$ws_GSM.Range("L73").Value = 15502.4
$ws_GSM.Range("M73").Value = -6096.875
$ws_GSM.Range("N73").Value = -17374.4
$ws_GSM.Range("H126").Value = 6878.5
$ws_GSM.Range("I126").Value = 3750
$ws_GSM.Range("J126").Value = 10007
$ws_GSM.Range("K126").Value = 11250
$ws_GSM.Range("L126").Value = 30021
$ws_GSM.Range("M126").Value = -8780
$ws_GSM.Range("N126").Value = -34961
$ws_GSM.Range("H132").Value = 1645020
$ws_GSM.Range("I132").Value = 1877306.6
$ws_GSM.Range("K132").Value = 5631919.800000001
$ws_GSM.Range("M132").Value = -5629389.800000001
$ws_GSM.Range("H134").Value = 18581.25
$ws_GSM.Range("J134").Value = 18581.25
$ws_GSM.Range("L134").Value = 55743.75
$ws_GSM.Range("N134").Value = -60813.75
$ws_GSM.Range("H136").Value = 51575.332
$ws_GSM.Range("J136").Value = 51575.332
$ws_GSM.Range("L136").Value = 154725.996
$ws_GSM.Range("N136").Value = -159825.996
# --- LTW ---
$ws_LTW.Range("H93").Value = 1151.375
$ws_LTW.Range("I93").Value = 1042.8334
$ws_LTW.Range("K93").Value = 1042.8334
$ws_LTW.Range("M93").Value = 205.1666
$ws_LTW.Range("H131").Value = 47999
$ws_LTW.Range("J131").Value = 47999
$ws_LTW.Range("L131").Value = 47999
$ws_LTW.Range("N131").Value = -58079
# --- WVR ---
$ws_WVR.Range("H96").Value = 2112.5
$ws_WVR.Range("J96").Value = 1650
$ws_WVR.Range("L96").Value = 1650
$ws_WVR.Range("N96").Value = -4396
$ws_WVR.Range("H107").Value = 470.35294
$ws_WVR.Range("I107").Value = 458.8889
$ws_WVR.Range("J107").Value = 483.25
$ws_WVR.Range("K107").Value = 1376.6667
$ws_WVR.Range("L107").Value = 1449.75
$ws_WVR.Range("M107").Value = 543.3333
$ws_WVR.Range("N107").Value = -5289.75
$ws_WVR.Range("H122").Value = 6858.4614
$ws_WVR.Range("I122").Value = 1845.1111
$ws_WVR.Range("K122").Value = 5535.3333
$ws_WVR.Range("M122").Value = -3085.3333
$ws_WVR.Range("H136").Value = 1924.4117
$ws_WVR.Range("I136").Value = 1417.5
$ws_WVR.Range("K136").Value = 4252.5
$ws_WVR.Range("M136").Value = -1702.5
